$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.923.23"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.995.24"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "544.47"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "138.35"
$ws.Range("E6").Value = "  +5.16%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "2.989.69"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  +14.07%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "33.69"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "3.481.85"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "61.924.68"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "2.996.10"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "462.14"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "13.22"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "0.649"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "78.79"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "12.47"
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").Value = "1.99"
$ws.Range("E29").Value = "  +5.54%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "25.26"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "2.32"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "54.58"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "448.07"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "0.0388"
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("D40").Value = "2.923.13"
$ws.Range("E40").Value = "  -7.17%  "
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "8.03"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  +7.26%  "
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "114.48"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "0.0₃0495"
$ws.Range("E50").Value = "  +2.71%  "
$ws.Range("E51").Value = "  -3.80%  "
